# Add uncommon skill "Bottled Plague" — update the card-stat tables.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rarity Distribution table (Table24, F6:I11) / "Mine" + note columns (J/K) ---
# Common row (8)
$ws.Range("E8").Value = 11
$ws.Range("J8").Value = 20
$ws.Range("K8").Value = "OK"

# Uncommon row (9) -- Bottled Plague is Uncommon
$ws.Range("E9").Value = 29
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 9

# Rare row (10)
$ws.Range("E10").Value = 14
$ws.Range("K10").Value = 3

# Totals row (11) -- J11 is SUM(J7:J10), recalculates automatically
$ws.Range("K11").Value = 13

# --- Type Distribution table (Table242, F16:I20) / "Mine" + note columns (J/K) ---
# Attack row (17)
$ws.Range("J17").Value = 22

# Skill row (18) -- Bottled Plague is a Skill
$ws.Range("J18").Value = 28

# Power row (19)
$ws.Range("J19").Value = 13
$ws.Range("K19").Value = 1

# Update the selected cell to match the author's final cursor position.
$ws.Range("E10").Select()
